$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investments")

# Insert a new row at position 19 (shifting the existing GLD row down to 20)
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with OILK data
$ws.Range("A19").Value = "OILK"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 1

# Copy formatting from the row above (XLRE, row 18) to match the standard style
$ws.Range("A18:C18").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the Max Weight for GLD (now row 20) from 0.1 to 1
$ws.Range("C20").Value = 1

# Update the active selection to match the target state
$ws.Range("C20").Select()
